$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.Value = "27.597.38"
$cell = $ws.Cells.Item(2, 5)
$cell.Value = "  +0.43%  "
$cell = $ws.Cells.Item(3, 4)
$cell.Value = "1.648.10"
$cell = $ws.Cells.Item(3, 5)
$cell.Value = "  -0.51%  "
$cell = $ws.Cells.Item(4, 5)
$cell.Value = "  +0.07%  "
$cell = $ws.Cells.Item(5, 4)
$cell.Value = "'213.30"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.Value = "  -0.35%  "
$cell = $ws.Cells.Item(6, 5)
$cell.Value = "  +4.88%  "
$cell = $ws.Cells.Item(7, 5)
$cell.Value = "  +0.04%  "
$cell = $ws.Cells.Item(8, 4)
$cell.Value = "'23.55"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 5)
$cell.Value = "  -2.30%  "
$cell = $ws.Cells.Item(9, 5)
$cell.Value = "  -2.14%  "
$cell = $ws.Cells.Item(10, 5)
$cell.Value = "  -0.86%  "
$cell = $ws.Cells.Item(11, 4)
$cell.Value = "'0.0890"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 5)
$cell.Value = "  +1.56%  "
$cell = $ws.Cells.Item(12, 4)
$cell.Value = "1.882.16"
$cell = $ws.Cells.Item(12, 5)
$cell.Value = "  -0.50%  "
$cell = $ws.Cells.Item(13, 4)
$cell.Value = "1.638.80"
$cell = $ws.Cells.Item(13, 5)
$cell.Value = "  -1.00%  "
$cell = $ws.Cells.Item(14, 4)
$cell.Value = "'0.588"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.Value = "  +4.70%  "
$cell = $ws.Cells.Item(15, 4)
$cell.Value = "'4.04"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 5)
$cell.Value = "  -1.99%  "
$cell = $ws.Cells.Item(16, 4)
$cell.Value = "'64.52"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 5)
$cell.Value = "  -1.84%  "
$cell = $ws.Cells.Item(17, 4)
$cell.Value = "27.579.14"
$cell = $ws.Cells.Item(17, 5)
$cell.Value = "  +0.34%  "
$cell = $ws.Cells.Item(18, 4)
$cell.Value = "'232.13"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.Value = "  -3.11%  "
$cell = $ws.Cells.Item(19, 4)
$cell.Value = "0.0₃0724"
$cell = $ws.Cells.Item(19, 5)
$cell.Value = "  -0.66%  "
$cell = $ws.Cells.Item(20, 4)
$cell.Value = "'7.62"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 5)
$cell.Value = "  +0.66%  "
$cell = $ws.Cells.Item(21, 5)
$cell.Value = "  +0.03%  "
$cell = $ws.Cells.Item(22, 5)
$cell.Value = "  -1.94%  "
$cell = $ws.Cells.Item(23, 4)
$cell.Value = "'9.74"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.Value = "  +3.97%  "
$cell = $ws.Cells.Item(24, 4)
$cell.Value = "'1.99"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.Value = "  -2.70%  "
$cell = $ws.Cells.Item(25, 4)
$cell.Value = "'148.69"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 5)
$cell.Value = "  +2.16%  "
$cell = $ws.Cells.Item(26, 4)
$cell.Value = "'7.04"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 5)
$cell.Value = "  -1.88%  "
$cell = $ws.Cells.Item(27, 5)
$cell.Value = "  +1.83%  "
$cell = $ws.Cells.Item(28, 5)
$cell.Value = "  -0.04%  "
$cell = $ws.Cells.Item(29, 4)
$cell.Value = "'15.63"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 5)
$cell.Value = "  -3.61%  "
$cell = $ws.Cells.Item(30, 5)
$cell.Value = "  -0.64%  "
$cell = $ws.Cells.Item(31, 4)
$cell.Value = "'0.0488"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 5)
$cell.Value = "  -2.10%  "
$cell = $ws.Cells.Item(32, 5)
$cell.Value = "  -0.05%  "
$cell = $ws.Cells.Item(33, 5)
$cell.Value = "  +3.68%  "
$cell = $ws.Cells.Item(34, 4)
$cell.Value = "1.433.52"
$cell = $ws.Cells.Item(34, 5)
$cell.Value = "  -0.77%  "
$cell = $ws.Cells.Item(35, 5)
$cell.Value = "  +2.91%  "
$cell = $ws.Cells.Item(36, 5)
$cell.Value = "  -0.19%  "
$cell = $ws.Cells.Item(37, 4)
$cell.Value = "'0.575"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 5)
$cell.Value = "  +0.75%  "
$cell = $ws.Cells.Item(38, 5)
$cell.Value = "  -3.51%  "
$cell = $ws.Cells.Item(40, 5)
$cell.Value = "  -3.49%  "
$cell = $ws.Cells.Item(41, 5)
$cell.Value = "  +0.06%  "
$cell = $ws.Cells.Item(42, 4)
$cell.Value = "'0.817"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.Value = "  +3.20%  "
$cell = $ws.Cells.Item(43, 4)
$cell.Value = "'5.49"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.Value = "  +1.29%  "
$cell = $ws.Cells.Item(44, 4)
$cell.Value = "'2.24"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 5)
$cell.Value = "  +0.93%  "
$cell = $ws.Cells.Item(45, 4)
$cell.Value = "'65.42"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 5)
$cell.Value = "  -1.59%  "
$cell = $ws.Cells.Item(46, 4)
$cell.Value = "1.791.39"
$cell = $ws.Cells.Item(47, 5)
$cell.Value = "  -0.18%  "
$cell = $ws.Cells.Item(48, 4)
$cell.Value = "'87.90"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 5)
$cell.Value = "  -0.57%  "
$cell = $ws.Cells.Item(49, 4)
$cell.Value = "0.0₆0108"
$cell = $ws.Cells.Item(49, 5)
$cell.Value = "  +1.63%  "
$cell = $ws.Cells.Item(50, 5)
$cell.Value = "  -1.65%  "
$cell = $ws.Cells.Item(51, 4)
$cell.Value = "'7.79"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 5)
$cell.Value = "  +0.13%  "
